# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Granada" (row 90 and 91),
# pushing the existing rows 90-147 down to 92-149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 90 (same position,
# done twice so both new rows land at 90 and 91 in final order).
$ws.Rows.Item(90).Insert()
$ws.Rows.Item(90).Insert()

# New row 90
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = 44777
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100104
$ws.Range("H90").Value = "Frutos de pepita"
$ws.Range("I90").Value = 100104001
$ws.Range("J90").Value = "Granada"
$ws.Range("K90").Value = "Wonderfull"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 120
$ws.Range("N90").Value = 15000
$ws.Range("O90").Value = 15000
$ws.Range("P90").Value = 15000
$ws.Range("Q90").Value = "$/bandeja 10 kilos"
$ws.Range("R90").Value = "Provincia de Limarí"
$ws.Range("S90").Value = 1500
$ws.Range("T90").Value = 10

# New row 91
$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 44777
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100104
$ws.Range("H91").Value = "Frutos de pepita"
$ws.Range("I91").Value = 100104001
$ws.Range("J91").Value = "Granada"
$ws.Range("K91").Value = "Wonderfull"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 200
$ws.Range("N91").Value = 14000
$ws.Range("O91").Value = 14000
$ws.Range("P91").Value = 14000
$ws.Range("Q91").Value = "$/bandeja 10 kilos granel"
$ws.Range("R91").Value = "Provincia de Limarí"
$ws.Range("S91").Value = 1400
$ws.Range("T91").Value = 10
